# Add the new LEGO set row (10280 - "Bouquet de fleurs") as row 22,
# matching the plain-text layout already used by every other data row
# on the sheet (ID/piece-count columns are numeric-looking text, and
# the three site columns with no vendor link are present but empty).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

# A leading apostrophe forces Excel to store each value as literal text
# - this keeps purely numeric-looking strings like "10280"/"756" as text
#   (matching column A/C elsewhere on the sheet) instead of becoming
#   numbers, and it lets an "empty" cell be written as a real, typed,
#   empty text value (columns G/H/I) rather than being cleared outright.
$ws.Cells.Item($row, 1).Value  = "'10280"
$ws.Cells.Item($row, 2).Value  = "'Bouquet de fleurs"
$ws.Cells.Item($row, 3).Value  = "'756"
$ws.Cells.Item($row, 4).Value  = "'The Botanical Collection"
$ws.Cells.Item($row, 5).Value  = "'https://www.lego.com/cdn/cs/set/assets/blt53711dac56e01b36/10280_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Cells.Item($row, 6).Value  = "'https://www.lego.com/fr-fr/product/10280"
$ws.Cells.Item($row, 7).Value  = "'"
$ws.Cells.Item($row, 8).Value  = "'"
$ws.Cells.Item($row, 9).Value  = "'"
$ws.Cells.Item($row, 10).Value = "'https://www.avenuedelabrique.com/lego-creator/10280-bouquet-de-fleurs/p6299"

# Reset the whole new row back to the default "Normal" style in one shot
# so it doesn't keep the transient quote-prefix formatting that the
# apostrophe trick above applies - same unstyled look as the rest of the
# table (row 21, etc.).
$ws.Range("A22:J22").Style = "Normal"
